$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in homework/attendance marks (value 5) for several students/columns
$ws.Range("C7:F7").Value = 5
$ws.Range("G21").Value = 5
$ws.Range("C24:E24").Value = 5
$ws.Range("C25:F25").Value = 5
$ws.Range("C29:F29").Value = 5
$ws.Range("C30").Value = 5

# Scroll the frozen pane down so its visible top-left cell is C12
# (freeze boundary itself - rows 1-3 / cols A-B - stays the same).
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 3

# Move the active selection to D30
$ws.Range("D30").Select()
